$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear all existing content first (old sheet used rows 1-95)
$ws.Range("A1:C95").ClearContents()

$colA = @(
    'Cluster name',
    '3535 Opal Meadow Heights Aged Care Community Meadow Heights',
    '3535 Opal Meadow Heights Aged Care Community Meadow Heights',
    'Acquire BPO Southbank',
    'Acquire BPO Southbank',
    'Al Haj Halal Meats Glenroy',
    'Al Haj Halal Meats Glenroy',
    'Al-Taqwa College Truganina',
    'Al-Taqwa College Truganina',
    'Broadmeadows Medical Centre Broadmeadows',
    'CS Square Caroline Springs',
    'Can Panel Cambellfield',
    'Can Panel Campbellfield',
    'City of Hobsons Bay Community',
    'City of Hobsons Bay Community',
    'City of Moreland Community',
    'City of Moreland Community',
    'City of Wyndham Community',
    'Coles Aurora Village Epping',
    'Coles Aurora Village Epping',
    'Coles Barkly Square Brunswick August',
    'Coles Campbellfield Plaza Campbellfield',
    'Coles Campbellfield Plaza Campbellfield',
    'Coles Coburg North Village',
    'Coles Coburg North Village',
    'Coles Coburg North Village Aug',
    'Coles Greenvale Shopping Centre',
    'Coles Greenvale Shopping Centre',
    'Coles Pakenham Place Shopping Centre',
    'Coles Pakenham Place Shopping Centre',
    'Costco Wholesale Epping',
    'Costco Wholesale Epping',
    'DRC Laverton Automotive Repairs Laverton North',
    'Direct Freight Express Cambellfield',
    'Direct Freight Express Campbellield',
    'Ernst and Young 8 Exhibition Street Melbourne',
    'Fitzroy Community School Fitzroy North',
    'Fitzroy Community School Fitzroy North',
    'Glenroy West Primary School',
    'Glenroy West Primary School',
    'Hamilton Marino 236 Jasper Road McKinnon',
    'Hamilton Marino 236 Jasper Road McKinnon',
    'Health Care Providers Association South Melbourne',
    'Health Care Providers Association South Melbourne',
    'IGA Meadow Heights Shopping Centre Meadow Heights',
    'IGA Meadow Heights Shopping Centre Meadow Heights',
    'Ilim College Kiewa Campus Boys Secondary Dallas',
    'Ilim Leaning Sanctuary Glenroy',
    'Ilim Learning Sanctuary Glenroy',
    'Industrial Galvanizers Valmont Coatings Campbellfield',
    'Industrial Galvanizers Valmont Coatings Campbellfield',
    'KFC Fawkner',
    'KFC Fawkner',
    'Malvern Health and Fitness Clinic Malvern',
    'McDonalds Thomastown II',
    'McDonalds Thomastown II',
    'National Gallery of Victoria Melbourne',
    'Newport Gardens Early Years Centre Newport',
    'Oporto Coolaroo',
    'Oporto Coolaroo',
    'Panorama Construction Site Whitehorse Rd Box Hill',
    'Panorama Construction Site Whitehorse Rd Box Hill',
    'Richmond Quarter 261-271 Bridge Road Construction Site Richmond',
    'Serco Mill Park',
    'The Royal Children''s Hospital Melbourne Emergency Department Parkville Tier 1B',
    'The Royal Children''s Hospital Melbourne Emergency Department Parkville Tier 1B',
    'Tip Top Warehouse Dandenong',
    'Tip Top Warehouse Dandenong',
    'Unilodge College Square Student Accommodation 570 Lygon Street Carlton',
    'Unilodge College Square Student Accommodation 570 Lygon Street Carlton',
    'Werribee Mercy Hospistal Emergency Department',
    'Werribee Mercy Hospital Emergency Department',
    'Western Health Footscray Hospital Emergency Department',
    'Who is Bunker Spreckels Cafe Elwood',
    'Woolworths Greenvale Lakes Roxburgh Park',
    'Woolworths Greenvale Lakes Roxburgh Park',
    'Yara Childcare Centre Truganina',
    'Yarra Childcare Centre Truganina'
)
$colB = @(
    'Active cases',
    22,
    25,
    9,
    12,
    50,
    53,
    7,
    12,
    5,
    8,
    5,
    5,
    9,
    10,
    6,
    7,
    6,
    5,
    6,
    5,
    8,
    10,
    15,
    27,
    7,
    6,
    7,
    7,
    8,
    20,
    21,
    5,
    10,
    6,
    5,
    49,
    55,
    6,
    7,
    6,
    11,
    13,
    16,
    6,
    7,
    5,
    15,
    15,
    15,
    17,
    9,
    11,
    5,
    6,
    7,
    9,
    5,
    16,
    17,
    27,
    28,
    9,
    5,
    7,
    8,
    8,
    10,
    13,
    14,
    7,
    7,
    5,
    5,
    6,
    10,
    13,
    15
)
$colC = @(
    'Exist',
    'old',
    'new',
    'new',
    'old',
    'new',
    'old',
    'new',
    'old',
    'old',
    'new',
    'new',
    'old',
    'new',
    'old',
    'old',
    'new',
    'new',
    'old',
    'new',
    'old',
    'old',
    'new',
    'old',
    'new',
    'old',
    'old',
    'new',
    'old',
    'new',
    'new',
    'old',
    'new',
    'new',
    'old',
    'old',
    'new',
    'old',
    'new',
    'old',
    'old',
    'new',
    'new',
    'old',
    'old',
    'new',
    'old',
    'old',
    'new',
    'old',
    'new',
    'new',
    'old',
    'old',
    'old',
    'new',
    'new',
    'old',
    'new',
    'old',
    'new',
    'old',
    'new',
    'old',
    'old',
    'new',
    'new',
    'old',
    'old',
    'new',
    'old',
    'new',
    'new',
    'old',
    'new',
    'old',
    'new',
    'old'
)

for ($i = 0; $i -lt $colA.Length; $i++) {
    $r = $i + 1
    $ws.Cells.Item($r, 1).Value = $colA[$i]
    $ws.Cells.Item($r, 2).Value = $colB[$i]
    $ws.Cells.Item($r, 3).Value = $colC[$i]
}

